# Edit described by the commit:
#   1. Slide 16's table (graphicFrame shape #3) switches to a different
#      table style (tableStyleId {DCB8C901-73DD-4064-872A-B9E6D548262B}).
#   2. The presentation's theme ("Integral") is swapped out for the plain
#      default "Office Theme" colour set (the old "Office Theme" colours
#      that used to live in theme1.xml become the active theme; the
#      "Integral" colours move to theme1.xml). The two themes only differ
#      in their colour scheme (fonts/effects are identical), so we
#      reproduce the colour swap by writing the "Office Theme" RGB values
#      into the active ThemeColorScheme.

function BGR($r, $g, $b) {
    # COM RGB values are packed 0x00BBGGRR (Windows BGR order).
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------------
$s16 = $p.Slides.Item(16)
$tbl = $s16.Shapes.Item(3).Table
$tbl.ApplyStyle("{DCB8C901-73DD-4064-872A-B9E6D548262B}")

# --- 2. Theme colours: Integral -> Office Theme ---------------------------------
# PpThemeColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$tcs = $p.SlideMaster.Theme.ThemeColorScheme
$tcs.Item(1).RGB  = BGR 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = BGR 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = BGR 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = BGR 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = BGR 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = BGR 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = BGR 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = BGR 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = BGR 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = BGR 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = BGR 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = BGR 0x95 0x4F 0x72   # folHlink 954F72
